# Phase-2 Start Add Interaction Comment Check.
# Updates the ticket-reference column (Y) on the NFTRTickets sheet with the
# newly created interaction-comment-check ticket numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NFTRTickets")

function Set-TicketId($addr, $value) {
    $cell = $ws.Range($addr)
    # These ids look numeric, so Excel would normally coerce them to a
    # number. Force text storage, then drop back to the cell's original
    # (default) style so no stray number-format/style is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TicketId "Y2" "270720001027"
Set-TicketId "Y3" "270720001028"
Set-TicketId "Y4" "270720001029"
Set-TicketId "Y5" "270720001030"
